$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Com")

# Insert a new row for the "Gro" (Grocery) prototype above the existing
# "Hsp" row, shifting the rest of the commercial prototype table down by one.
[void]$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = "Gro"
$ws.Range("C8").Value = "Area-ft2-BA"
$ws.Range("D8").Value = 50000

# Leave the selection where the user's last entry landed.
[void]$ws.Range("D9").Select()
